# Insert a new bullet list item after the "LDAP authentication." paragraph,
# matching the style/formatting of the existing list items, and before
# the trailing _GoBack bookmark.

$d = $word.ActiveDocument

# Locate the paragraph whose text is "LDAP authentication." and remember
# its 1-based index within the Paragraphs collection.
$target = $null
$targetIndex = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text.TrimEnd("`r") -eq "LDAP authentication.") {
        $target = $p
        $targetIndex = $i
        break
    }
}

if ($target -eq $null) {
    throw "Could not find 'LDAP authentication.' paragraph"
}

# Insert a new paragraph right after this one (inherits the same
# ListParagraph style / numbering / run formatting).
$target.Range.InsertParagraphAfter() | Out-Null

# Re-fetch the freshly inserted paragraph by index (object references
# captured before the structural edit are no longer valid).
$newPara = $d.Paragraphs.Item($targetIndex + 1)

# Fill in the new bullet text.
$newPara.Range.Text = "Add functionality to add aliquots by scanning the barcode or CSV file into the layout."
